$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column C (customer_code_raw) to Text so numeric-looking codes are preserved as strings
$ws.Range("C567:C592").NumberFormat = "@"

$rows = @(
  @{ Row=567; A=46001; B='한만교'; C='10108'; D='수'; E=12; F=10; G='셔2'; H=''; I=''; J=701; K=13 },
  @{ Row=568; A=46002; B='김주영'; C='10082'; D='목'; E=12; F=11; G='상1,하1'; H=''; I=''; J=701; K=17 },
  @{ Row=569; A=46003; B='조규종'; C='6098'; D='금'; E=12; F=12; G='상1,하1,셔3'; H=''; I=''; J=701; K=21 },
  @{ Row=570; A=46000; B='김정대'; C='10079'; D='화'; E=12; F=9; G='상1,중1,하2'; H=''; I=''; J=702; K=9 },
  @{ Row=571; A=46002; B='김영식'; C='9905'; D='목'; E=12; F=11; G='상1,하2,셔3'; H=''; I=''; J=702; K=17 },
  @{ Row=572; A=46003; B='이주환'; C='2953'; D='금'; E=12; F=12; G='상1,중1,하1,셔1'; H=''; I=''; J=702; K=21 },
  @{ Row=573; A=46000; B='김동수'; C='9491'; D='화'; E=12; F=9; G='하1 수선'; H=''; I=''; J=703; K=9 },
  @{ Row=574; A=46003; B='신덕호'; C='7312'; D='금'; E=12; F=12; G='코2,하1,셔2'; H=''; I=''; J=703; K=21 },
  @{ Row=575; A=46003; B='강두현'; C='10112'; D='금'; E=12; F=12; G='대여복1'; H=''; I=''; J=704; K=21 },
  @{ Row=576; A=46006; B='권재운'; C='10094'; D='월'; E=12; F=15; G='상1,하1,셔3'; H=''; I=''; J=711; K=5 },
  @{ Row=577; A=46007; B='박인영'; C='10107'; D='화'; E=12; F=16; G='상1,하1,셔2'; H=''; I=''; J=711; K=9 },
  @{ Row=578; A=46008; B='임요셉'; C='5115'; D='수'; E=12; F=17; G='셔2'; H=''; I='택배'; J=711; K=13 },
  @{ Row=579; A=46013; B='김민수'; C='5582'; D='월'; E=12; F=22; G='코트수선'; H=''; I=''; J=721; K=5 },
  @{ Row=580; A=46014; B='김유현'; C='7464'; D='화'; E=12; F=23; G='상1,하2,셔2'; H=''; I=''; J=721; K=9 },
  @{ Row=581; A=46018; B='성일용'; C='3811'; D='토'; E=12; F=27; G='상1,하1,셔1'; H=''; I=''; J=721; K=25 },
  @{ Row=582; A=46013; B='장창석'; C='10109'; D='월'; E=12; F=22; G='상1,하1,셔2'; H=''; I=''; J=722; K=5 },
  @{ Row=583; A=46014; B='이세한'; C='7399'; D='화'; E=12; F=23; G='상1,하1,셔2'; H=''; I=''; J=722; K=9 },
  @{ Row=584; A=46013; B='염기태'; C='10110'; D='월'; E=12; F=22; G='상1,하2,셔2'; H=''; I=''; J=723; K=5 },
  @{ Row=585; A=46014; B='배주원'; C='10105'; D='화'; E=12; F=23; G='상2,하2,셔3'; H=''; I=''; J=723; K=9 },
  @{ Row=586; A=46014; B='이현철'; C='10102'; D='화'; E=12; F=23; G='상1,하1,셔1'; H=''; I=''; J=724; K=9 },
  @{ Row=587; A=46014; B='이상윤'; C='10114'; D='화'; E=12; F=23; G='코1,상1,하1,셔1'; H=''; I=''; J=725; K=9 },
  @{ Row=588; A=46014; B='김재현'; C='6053'; D='화'; E=12; F=23; G='코1'; H=''; I=''; J=726; K=9 },
  @{ Row=589; A=46014; B='김준후'; C='9923'; D='화'; E=12; F=23; G='상1,하1,셔2'; H=''; I=''; J=727; K=9 },
  @{ Row=590; A=46020; B='이용정'; C='10111'; D='월'; E=12; F=29; G='대여복2'; H=''; I=''; J=731; K=5 },
  @{ Row=591; A=46021; B='박정환'; C='10113'; D='화'; E=12; F=30; G='코1,셔1'; H=''; I=''; J=731; K=9 },
  @{ Row=592; A=46020; B='박상언'; C='10106'; D='월'; E=12; F=29; G='상1,하2,셔2'; H=''; I=''; J=732; K=5 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    if ($r.H -ne "") { $ws.Cells.Item($row, 8).Value = $r.H }
    if ($r.I -ne "") { $ws.Cells.Item($row, 9).Value = $r.I }
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
}

# Column A carries the existing date-style (style index inherited from header rows); re-apply it explicitly
$ws.Range("A567:A592").NumberFormat = $ws.Range("A566").NumberFormat

Write-Host "Updated rows 567-592"
